$wb = $excel.ActiveWorkbook

# --- SS_att sheet: revised summary statistics (rows 2 & 3) ---
$ssAtt = $wb.Worksheets.Item("SS_att")

$ssAtt.Range("B2").Value = 32.333333333333336
$ssAtt.Range("C2").Value = 34.15
$ssAtt.Range("D2").Value = 37.134328358208954
$ssAtt.Range("E2").Value = 34.752808988764045
$ssAtt.Range("L2").Value = 0.572023905328648

$ssAtt.Range("B3").Value = 3.3381533778207677
$ssAtt.Range("C3").Value = 2.9165246131670202
$ssAtt.Range("D3").Value = 3.1781063232804376
$ssAtt.Range("E3").Value = 1.8159027540097266

# --- Attrition sheet: fix formulas that referenced the wrong column (E -> D) ---
$attrition = $wb.Worksheets.Item("Attrition")

$attrition.Range("D6").Formula = "=ROUND(SS_att!D2,0)"
$attrition.Range("D7").Formula = "=CONCATENATE(""("",ROUND(SS_att!D3,1),"")"")"
$attrition.Range("D10").Formula = "=ROUND(SS_att!D10,2)"
$attrition.Range("D11").Formula = "=CONCATENATE(""("",ROUND(SS_att!D11,2),"")"")"

$wb.Application.CalculateFull()
